$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.920330882072449
$ws.Range("B1").Value = 1.92260217666626
$ws.Range("C1").Value = 7.92901611328125
$ws.Range("D1").Value = 0.928367018699646
$ws.Range("E1").Value = 0.4212445020675659
